$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the "not found" placeholder row for the SBI statement
# (it previously held the full "Attachment Found" record, which now
# moves down to row 3). Rename the subject while we're at it.
$ws.Range("A2").Value = "SBI Statement"
$ws.Range("B2").Value = "Attachment Not Found"
$ws.Range("C2:I2").ClearContents()

# Row 3 used to be the short "Attachment Not Found" row for SBI; it now
# carries the full record (renamed subject + refreshed 06-24-2022 path).
$ws.Range("A3").Value = "SBI Statement"
$ws.Range("B3").Value = "Attachment Found"
$ws.Range("C3").Value = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\06-24-2022\SBI Bank\Vidya Sagar Reddy.pdf"
$ws.Range("D3").Value = "Password Found"
$ws.Range("E3").Value = "Password Removed"
$ws.Range("F3").Value = "1182"
$ws.Range("G3").Value = "280701501966,"
$ws.Range("H3").Value = "U72200TG2014PTC092878"
$ws.Range("I3").Value = "May-2022"

# Rows 4-6: rename subjects and refresh the attachment paths to
# 06-24-2022 (remaining columns were already correct).
$ws.Range("A4").Value = "AXIS Statement"
$ws.Range("C4").Value = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\06-24-2022\Axis Bank\Vidya Sagar Reddy.pdf"

$ws.Range("A5").Value = "ICICI Statement"
$ws.Range("C5").Value = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\06-24-2022\ICICI Bank\Vidya Sagar Reddy.pdf"

$ws.Range("A6").Value = "HDFC Statement"
$ws.Range("C6").Value = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\06-24-2022\HDFC Bank\Vidya Sagar Reddy.pdf"

# View tidy-up: drop the old far-right J-column selection and scroll
# the sheet back towards the left (closest available to topLeftCell=B1).
$ws.Range("A1").Select()

# Column A no longer needs to be as wide since the new subject labels
# are shorter; column I's width entry is no longer used.
$ws.Columns.Item(1).ColumnWidth = 13.6
